$wb = $excel.ActiveWorkbook

# --- Sheet "List" (book detail row): accountId -> owner (printed), title -> bookTitle ---
$ws1 = $wb.Worksheets.Item("List")

$ws1.Range("B1").Value = "`${msg.getProperty('book_owner')}"
$ws1.Range("B2").Value = "`${printer.print(book.owner)}"
$ws1.Range("C1").Value = "`${msg.getProperty('book_bookTitle')}"
$ws1.Range("C2").Value = "`${book.bookTitle}"

# --- Sheet "Search": same rename of criteria fields, plus split of the old combined row ---
$ws2 = $wb.Worksheets.Item("Search")

# Row 5 used to be the book_title/title search row -> now book_owner/owner
$ws2.Range("A5").Value = "`${msg.getProperty('book_owner')}"
$ws2.Range("B5").Value = "`${owner}"

# Row 6 used to hold numberOfPages + range_from/range_to all together.
# It is now split: row 6 keeps only the (renamed) bookTitle row, and a new
# row 7 carries the numberOfPages / range_from / range_to fields that used
# to live on row 6.
$ws2.Range("A6").Value = "`${msg.getProperty('book_bookTitle')}"
$ws2.Range("B6").Value = "`${bookTitle}"
$ws2.Range("C6").Value = ""
$ws2.Range("D6").Value = ""
$ws2.Range("E6").Value = ""

$ws2.Range("A7").Value = "`${msg.getProperty('book_numberOfPages')}"
$ws2.Range("B7").Value = "`${msg.getProperty('range_from')}"
$ws2.Range("C7").Value = "`${numberOfPagesFrom}"
$ws2.Range("D7").Value = "`${msg.getProperty('range_to')}"
$ws2.Range("E7").Value = "`${numberOfPagesTo}"
